$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 3.173141015532005
$ws.Cells.Item(2, 4).Value = 4.035425324497355
$ws.Cells.Item(2, 5).Value = 16.53723457708421
$ws.Cells.Item(2, 6).Value = 21.89869758807772
$ws.Cells.Item(2, 7).Value = 3.605196027482979
$ws.Cells.Item(2, 9).Value = 18.31535829197669
$ws.Cells.Item(2, 11).Value = 13.11107939293262
$ws.Cells.Item(2, 14).Value = 16.18344107474968
$ws.Cells.Item(2, 15).Value = 19.39020053596778

$ws.Cells.Item(3, 3).Value = 3.061578761134364
$ws.Cells.Item(3, 4).Value = 4.001823076210073
$ws.Cells.Item(3, 5).Value = 15.59608990512626
$ws.Cells.Item(3, 6).Value = 21.83568393080377
$ws.Cells.Item(3, 7).Value = 3.607469045620589
$ws.Cells.Item(3, 9).Value = 18.31410290542723
$ws.Cells.Item(3, 11).Value = 12.53672800248722
$ws.Cells.Item(3, 14).Value = 16.21548897535106
$ws.Cells.Item(3, 15).Value = 19.41247630692527

$ws.Cells.Item(4, 3).Value = 2.990320470498073
$ws.Cells.Item(4, 4).Value = 3.980827431120664
$ws.Cells.Item(4, 5).Value = 14.99323847450719
$ws.Cells.Item(4, 6).Value = 21.80501445857942
$ws.Cells.Item(4, 7).Value = 3.608937064847609
$ws.Cells.Item(4, 9).Value = 18.31930158091502
$ws.Cells.Item(4, 11).Value = 12.17133699759901
$ws.Cells.Item(4, 14).Value = 16.23697042658797
$ws.Cells.Item(4, 15).Value = 19.43246943495614

$ws.Cells.Item(5, 3).Value = 2.96061425917886
$ws.Cells.Item(5, 4).Value = 3.972184185711812
$ws.Cells.Item(5, 5).Value = 14.74155777864946
$ws.Cells.Item(5, 6).Value = 21.79453942385345
$ws.Cells.Item(5, 7).Value = 3.60955355724154
$ws.Cells.Item(5, 9).Value = 18.32291771466971
$ws.Cells.Item(5, 11).Value = 12.01943549805005
$ws.Cells.Item(5, 14).Value = 16.24617861553001
$ws.Cells.Item(5, 15).Value = 19.4421980334461

$ws.Cells.Item(6, 3).Value = 2.955642039611275
$ws.Cells.Item(6, 4).Value = 3.970743826655564
$ws.Cells.Item(6, 5).Value = 14.69941207327079
$ws.Cells.Item(6, 6).Value = 21.79292235685353
$ws.Cells.Item(6, 7).Value = 3.609657030200707
$ws.Cells.Item(6, 9).Value = 18.32360849952544
$ws.Cells.Item(6, 11).Value = 11.99403749572
$ws.Cells.Item(6, 14).Value = 16.24773508941716
$ws.Cells.Item(6, 15).Value = 19.4439087423557

$ws.Cells.Item(7, 3).Value = 2.989922510101493
$ws.Cells.Item(7, 4).Value = 3.9807112135251
$ws.Cells.Item(7, 5).Value = 14.9898681762325
$ws.Cells.Item(7, 6).Value = 21.80486499173729
$ws.Cells.Item(7, 7).Value = 3.608945305051557
$ws.Cells.Item(7, 9).Value = 18.31934429103361
$ws.Cells.Item(7, 11).Value = 12.16930026863308
$ws.Cells.Item(7, 14).Value = 16.23709277097466
$ws.Cells.Item(7, 15).Value = 19.4325942457728

$ws.Cells.Item(8, 3).Value = 3.135261467943206
$ws.Cells.Item(8, 4).Value = 4.023917212782567
$ws.Cells.Item(8, 5).Value = 16.2180590568508
$ws.Cells.Item(8, 6).Value = 21.87531019364223
$ws.Cells.Item(8, 7).Value = 3.605964780244687
$ws.Cells.Item(8, 9).Value = 18.31368476134472
$ws.Cells.Item(8, 11).Value = 12.91579888210009
$ws.Cells.Item(8, 14).Value = 16.19411721590336
$ws.Cells.Item(8, 15).Value = 19.39656652876989

$ws.Cells.Item(9, 3).Value = 3.39739539030199
$ws.Cells.Item(9, 4).Value = 4.105580841975652
$ws.Cells.Item(9, 5).Value = 18.51872983103294
$ws.Cells.Item(9, 6).Value = 22.07669464897036
$ws.Cells.Item(9, 7).Value = 3.600691384694117
$ws.Cells.Item(9, 9).Value = 18.35005030638882
$ws.Cells.Item(9, 11).Value = 14.27110911920816
$ws.Cells.Item(9, 14).Value = 16.12412472457299
$ws.Cells.Item(9, 15).Value = 19.37629391610841

$ws.Cells.Item(10, 3).Value = 3.57490649208226
$ws.Cells.Item(10, 4).Value = 4.163479825091446
$ws.Cells.Item(10, 5).Value = 20.16701408909895
$ws.Cells.Item(10, 6).Value = 22.26246997376724
$ws.Cells.Item(10, 7).Value = 3.597161335318237
$ws.Cells.Item(10, 9).Value = 18.40575700139582
$ws.Cells.Item(10, 11).Value = 15.19226507805495
$ws.Cells.Item(10, 14).Value = 16.08136772438491
$ws.Cells.Item(10, 15).Value = 19.39240895822138

$ws.Cells.Item(11, 3).Value = 3.652195308303011
$ws.Cells.Item(11, 4).Value = 4.189315881294052
$ws.Cells.Item(11, 5).Value = 20.87461755499322
$ws.Cells.Item(11, 6).Value = 22.35497958516645
$ws.Cells.Item(11, 7).Value = 3.595629329446449
$ws.Cells.Item(11, 9).Value = 18.43737801442916
$ws.Cells.Item(11, 11).Value = 15.59376333187598
$ws.Cells.Item(11, 14).Value = 16.06378993341679
$ws.Cells.Item(11, 15).Value = 19.40651472592361

$ws.Cells.Item(12, 3).Value = 3.68095152724795
$ws.Cells.Item(12, 4).Value = 4.199023115053564
$ws.Cells.Item(12, 5).Value = 21.1365358659214
$ws.Cells.Item(12, 6).Value = 22.39113852188376
$ws.Cells.Item(12, 7).Value = 3.595059750027636
$ws.Cells.Item(12, 9).Value = 18.45025147543688
$ws.Cells.Item(12, 11).Value = 15.74318731715966
$ws.Cells.Item(12, 14).Value = 16.05740229262519
$ws.Cells.Item(12, 15).Value = 19.41283192653635

$ws.Cells.Item(13, 3).Value = 3.674781300927937
$ws.Cells.Item(13, 4).Value = 4.196935946637916
$ws.Cells.Item(13, 5).Value = 21.08039497642221
$ws.Cells.Item(13, 6).Value = 22.38330128880224
$ws.Cells.Item(13, 7).Value = 3.595181950512909
$ws.Cells.Item(13, 9).Value = 18.44743902420678
$ws.Cells.Item(13, 11).Value = 15.71112379285945
$ws.Cells.Item(13, 14).Value = 16.05876604471556
$ws.Cells.Item(13, 15).Value = 19.41142800180721

$ws.Cells.Item(14, 3).Value = 3.654571419606858
$ws.Cells.Item(14, 4).Value = 4.190116053770961
$ws.Cells.Item(14, 5).Value = 20.8962865044877
$ws.Cells.Item(14, 6).Value = 22.35793194413186
$ws.Cells.Item(14, 7).Value = 3.595582258543929
$ws.Cells.Item(14, 9).Value = 18.43841912543974
$ws.Cells.Item(14, 11).Value = 15.60610941230807
$ws.Cells.Item(14, 14).Value = 16.06325903682006
$ws.Cells.Item(14, 15).Value = 19.40701488833004

$ws.Cells.Item(15, 3).Value = 3.642125316508765
$ws.Cells.Item(15, 4).Value = 4.185928614326929
$ws.Cells.Item(15, 5).Value = 20.78272971649286
$ws.Cells.Item(15, 6).Value = 22.3425386101276
$ws.Cells.Item(15, 7).Value = 3.59582883198664
$ws.Cells.Item(15, 9).Value = 18.43301115386425
$ws.Cells.Item(15, 11).Value = 15.54144207227748
$ws.Cells.Item(15, 14).Value = 16.06604609787149
$ws.Cells.Item(15, 15).Value = 19.40443880978365

$ws.Cells.Item(16, 3).Value = 3.569784621606364
$ws.Cells.Item(16, 4).Value = 4.16178091736611
$ws.Cells.Item(16, 5).Value = 20.1199240807476
$ws.Cells.Item(16, 6).Value = 22.2565833358173
$ws.Cells.Item(16, 7).Value = 3.59726293613773
$ws.Cells.Item(16, 9).Value = 18.40381656901016
$ws.Cells.Item(16, 11).Value = 15.16566502368815
$ws.Cells.Item(16, 14).Value = 16.08255410510793
$ws.Cells.Item(16, 15).Value = 19.39162359092817

$ws.Cells.Item(17, 3).Value = 3.524509020299491
$ws.Cells.Item(17, 4).Value = 4.146835469670913
$ws.Cells.Item(17, 5).Value = 19.7025276537457
$ws.Cells.Item(17, 6).Value = 22.20588604171683
$ws.Cells.Item(17, 7).Value = 3.598161580539632
$ws.Cells.Item(17, 9).Value = 18.38751279727215
$ws.Cells.Item(17, 11).Value = 14.93057489040494
$ws.Cells.Item(17, 14).Value = 16.09316042585988
$ws.Cells.Item(17, 15).Value = 19.3854984819811

$ws.Cells.Item(18, 3).Value = 3.498142527446154
$ws.Cells.Item(18, 4).Value = 4.13819218492818
$ws.Cells.Item(18, 5).Value = 19.45847731759811
$ws.Cells.Item(18, 6).Value = 22.17748077836219
$ws.Cells.Item(18, 7).Value = 3.598685409879681
$ws.Cells.Item(18, 9).Value = 18.37872669087093
$ws.Cells.Item(18, 11).Value = 14.79371078959373
$ws.Cells.Item(18, 14).Value = 16.09943719204016
$ws.Cells.Item(18, 15).Value = 19.38261316086647

$ws.Cells.Item(19, 3).Value = 3.489159883207479
$ws.Cells.Item(19, 4).Value = 4.135257760216201
$ws.Cells.Item(19, 5).Value = 19.3751623292219
$ws.Cells.Item(19, 6).Value = 22.16799347838312
$ws.Cells.Item(19, 7).Value = 3.598863965452648
$ws.Cells.Item(19, 9).Value = 18.37585353038371
$ws.Cells.Item(19, 11).Value = 14.74709123821055
$ws.Cells.Item(19, 14).Value = 16.10159269341818
$ws.Cells.Item(19, 15).Value = 19.38174570031438

$ws.Cells.Item(20, 3).Value = 3.529362442878344
$ws.Cells.Item(20, 4).Value = 4.14843133818184
$ws.Cells.Item(20, 5).Value = 19.74737136243604
$ws.Cells.Item(20, 6).Value = 22.2112049261536
$ws.Cells.Item(20, 7).Value = 3.598065199142189
$ws.Cells.Item(20, 9).Value = 18.3891871775585
$ws.Cells.Item(20, 11).Value = 14.95577175728344
$ws.Cells.Item(20, 14).Value = 16.09201312370917
$ws.Cells.Item(20, 15).Value = 19.38608449873916

$ws.Cells.Item(21, 3).Value = 3.660521529181337
$ws.Cells.Item(21, 4).Value = 4.192121324848939
$ws.Cells.Item(21, 5).Value = 20.95052716553865
$ws.Cells.Item(21, 6).Value = 22.36535312802666
$ws.Cells.Item(21, 7).Value = 3.59546439229878
$ws.Cells.Item(21, 9).Value = 18.44104411907471
$ws.Cells.Item(21, 11).Value = 15.63702630681225
$ws.Cells.Item(21, 14).Value = 16.06193204868809
$ws.Cells.Item(21, 15).Value = 19.40828464229974

$ws.Cells.Item(22, 3).Value = 3.743255853002812
$ws.Cells.Item(22, 4).Value = 4.220228587482958
$ws.Cells.Item(22, 5).Value = 21.70169354592673
$ws.Cells.Item(22, 6).Value = 22.47265783600159
$ws.Cells.Item(22, 7).Value = 3.593826127434578
$ws.Cells.Item(22, 9).Value = 18.48017456251313
$ws.Cells.Item(22, 11).Value = 16.06699227994576
$ws.Cells.Item(22, 14).Value = 16.04383817319593
$ws.Cells.Item(22, 15).Value = 19.42848039552975

$ws.Cells.Item(23, 3).Value = 3.699376189602931
$ws.Cells.Item(23, 4).Value = 4.205269400869952
$ws.Cells.Item(23, 5).Value = 21.30398841869939
$ws.Cells.Item(23, 6).Value = 22.41479522226064
$ws.Cells.Item(23, 7).Value = 3.594694891177742
$ws.Cells.Item(23, 9).Value = 18.45881213587053
$ws.Cells.Item(23, 11).Value = 15.83893565675368
$ws.Cells.Item(23, 14).Value = 16.05335213273939
$ws.Cells.Item(23, 15).Value = 19.41718104555872

$ws.Cells.Item(24, 3).Value = 3.527169259000038
$ws.Cells.Item(24, 4).Value = 4.147710004378354
$ws.Cells.Item(24, 5).Value = 19.72711024448058
$ws.Cells.Item(24, 6).Value = 22.20879794763623
$ws.Cells.Item(24, 7).Value = 3.598108750780727
$ws.Cells.Item(24, 9).Value = 18.38842836094032
$ws.Cells.Item(24, 11).Value = 14.94438556553282
$ws.Cells.Item(24, 14).Value = 16.09253126119653
$ws.Cells.Item(24, 15).Value = 19.38581757920419

$ws.Cells.Item(25, 3).Value = 3.329056837437599
$ws.Cells.Item(25, 4).Value = 4.083840786152327
$ws.Cells.Item(25, 5).Value = 17.87407867818448
$ws.Cells.Item(25, 6).Value = 22.01550891827736
$ws.Cells.Item(25, 7).Value = 3.602057222888542
$ws.Cells.Item(25, 9).Value = 18.33512035540051
$ws.Cells.Item(25, 11).Value = 13.9170364687107
$ws.Cells.Item(25, 14).Value = 16.14153472852346
$ws.Cells.Item(25, 15).Value = 19.37634994729715
